# Update the expense rows: rename some categories, change amounts/dates,
# and remove the last two rows (Travel / Rent) entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Gym Monthly"
$ws.Range("B2").Value = 800
$ws.Range("C2").Value = "2025-11-30T00:00:00.000Z"

$ws.Range("A3").Value = "New Clothes"
$ws.Range("B3").Value = 500
$ws.Range("C3").Value = "2025-11-20T00:00:00.000Z"

$ws.Range("A4").Value = "Groceries"
$ws.Range("B4").Value = 2000
$ws.Range("C4").Value = "2025-11-16T00:00:00.000Z"

$ws.Range("A5").Value = "House Rent"
$ws.Range("B5").Value = 5000
$ws.Range("C5").Value = "2025-11-15T00:00:00.000Z"

# Rows 6 and 7 (old "Travel" / "Rent" entries) no longer exist; remove them.
$ws.Rows("6:7").Delete()
